$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Insert a new row above row 15 on Sheet1 (pushes old rows 15-47 down to 16-48)
# and fill it in with the new "New 200m transport / depth" case.
$ws1.Rows(15).Insert()

$ws1.Range("A15").Value = "New 200m transport / depth"
$ws1.Range("B15").Value = -2.4333999999999998
$ws1.Range("C15").Value = 0.39710000000000001
$ws1.Range("D15").Value = -1.4021999999999999
$ws1.Range("E15").Value = 0.90859999999999996
$ws1.Range("F15").Value = -0.23330000000000001
$ws1.Range("G15").Value = -0.3029
$ws1.Range("H15").Value = -0.14630000000000001
$ws1.Range("I15").Value = 0.070800000000000002

# Make Sheet1 the active sheet/tab, with the current selection on H16
# (and Sheet2 is no longer the active tab).
$ws1.Select()
$ws1.Range("H16").Select()
